$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptos list: updated Price (col D) and Volume(1h) (col E) values.
# D values that look like plain decimal numbers are entered with a leading
# apostrophe so Excel keeps them as text (matching the source data, which
# stores every price/volume as a text string, not a number).
$ws.Range("D2").Value = "65.976.91"
$ws.Range("E2").Value = "  +6.80%  "
$ws.Range("D3").Value = "3.011.93"
$ws.Range("E3").Value = "  +4.26%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'582.62"
$ws.Range("E5").Value = "  +2.83%  "
$ws.Range("D6").Value = "'161.53"
$ws.Range("E6").Value = "  +13.00%  "
$ws.Range("D8").Value = "3.008.58"
$ws.Range("E8").Value = "  +4.27%  "
$ws.Range("E9").Value = "  +3.13%  "
$ws.Range("D10").Value = "'6.98"
$ws.Range("E10").Value = "  +0.29%  "
$ws.Range("E11").Value = "  +6.33%  "
$ws.Range("E12").Value = "  +5.81%  "
$ws.Range("E13").Value = "  +8.93%  "
$ws.Range("E14").Value = "  +8.94%  "
$ws.Range("E15").Value = "  +0.85%  "
$ws.Range("D16").Value = "65.986.35"
$ws.Range("E16").Value = "  +6.91%  "
$ws.Range("D17").Value = "3.512.86"
$ws.Range("E17").Value = "  +4.29%  "
$ws.Range("E18").Value = "  +6.37%  "
$ws.Range("D19").Value = "3.011.90"
$ws.Range("E19").Value = "  +4.16%  "
$ws.Range("D20").Value = "'457.38"
$ws.Range("E20").Value = "  +6.68%  "
$ws.Range("D21").Value = "'13.92"
$ws.Range("E21").Value = "  +7.28%  "
$ws.Range("E22").Value = "  +4.79%  "
$ws.Range("D23").Value = "'7.31"
$ws.Range("E23").Value = "  +6.66%  "
$ws.Range("D24").Value = "'82.23"
$ws.Range("E24").Value = "  +4.21%  "
$ws.Range("D25").Value = "'2.29"
$ws.Range("E25").Value = "  +13.12%  "
$ws.Range("D26").Value = "'12.35"
$ws.Range("E26").Value = "  +2.62%  "
$ws.Range("D27").Value = "'10.63"
$ws.Range("E27").Value = "  +5.85%  "
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("D29").Value = "'8.13"
$ws.Range("E29").Value = "  +17.08%  "
$ws.Range("D30").Value = "'2.34"
$ws.Range("E30").Value = "  +16.24%  "
$ws.Range("E31").Value = "  -5.48%  "
$ws.Range("E32").Value = "  +4.27%  "
$ws.Range("D33").Value = "'26.94"
$ws.Range("E33").Value = "  +5.31%  "
$ws.Range("E34").Value = "  +2.95%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("D36").Value = "'0.993"
$ws.Range("E36").Value = "  +4.73%  "
$ws.Range("D37").Value = "'5.77"
$ws.Range("E37").Value = "  +7.46%  "
$ws.Range("E38").Value = "  +12.24%  "
$ws.Range("D39").Value = "'49.92"
$ws.Range("E39").Value = "  +2.24%  "
$ws.Range("E40").Value = "  +6.46%  "
$ws.Range("E41").Value = "  +13.28%  "
$ws.Range("E42").Value = "  +5.84%  "
$ws.Range("E43").Value = "  +10.35%  "
$ws.Range("D44").Value = "'8.47"
$ws.Range("E44").Value = "  +3.94%  "
$ws.Range("D45").Value = "'383.93"
$ws.Range("E45").Value = "  +11.47%  "
$ws.Range("E46").Value = "  +5.89%  "
$ws.Range("D47").Value = "2.786.09"
$ws.Range("E47").Value = "  +3.63%  "
$ws.Range("D48").Value = "'134.73"
$ws.Range("E48").Value = "  +2.59%  "
$ws.Range("D50").Value = "'23.82"
$ws.Range("E50").Value = "  +10.59%  "
$ws.Range("D51").Value = "'0.106"
$ws.Range("E51").Value = "  +3.95%  "
